$d = $word.ActiveDocument
$p1 = $d.Paragraphs.Item(1)
$r = $p1.Range
$r.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item(2)
$ir = $p2.Range
$ir.InsertAfter("#2nd changes")

# Now find "nd" within paragraph 2 range and set superscript
$p2r = $p2.Range
Write-Output ("p2 start=" + $p2r.Start + " end=" + $p2r.End)

$findRange = $d.Range($p2r.Start, $p2r.End)
$findRange.Find.Execute("nd", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output ("Found: " + $findRange.Find.Found)
Write-Output ("Found range: [" + $findRange.Text + "] start=" + $findRange.Start + " end=" + $findRange.End)
$findRange.Font.Superscript = $true

Write-Output ("p2 text after: [" + $p2.Range.Text + "]")
